# "depployment first try finished"
# Update the workDescription column header to clarify the format, and
# drop the trailing period from the first row's description so it reads
# consistently with the "separated by ." convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("E1").Value = "workDescription(separated by .)"
$ws.Range("E2").Value = "Be responsible for product test data process system development.Product line software maintenance and program development"

# Scroll the window so column E is in view and leave the selection on E14,
# matching the author's final on-screen state when they saved the file.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("E14").Select()
